$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2100591715976331
$ws.Range("C2").Value = 0.5266272189349113
$ws.Range("J2").Value = 0.02366863905325444
$ws.Range("P2").Value = 0.1420118343195266
$ws.Range("S2").Value = 0.09763313609467456
$ws.Range("B3").Value = 0.005434782608695652
$ws.Range("C3").Value = 0.03804347826086957
$ws.Range("J3").Value = 0.02717391304347826
$ws.Range("P3").Value = 0.7065217391304348
$ws.Range("S3").Value = 0.2228260869565217
$ws.Range("J4").Value = 0.01923076923076923
$ws.Range("P4").Value = 0.6538461538461539
$ws.Range("S4").Value = 0.3269230769230769
$ws.Range("B6").Value = 0.03643724696356275
$ws.Range("D6").Value = 0.01619433198380567
$ws.Range("F6").Value = 0.06477732793522267
$ws.Range("J6").Value = 0.2631578947368421
$ws.Range("O6").Value = 0.03643724696356275
$ws.Range("Q6").Value = 0.145748987854251
$ws.Range("R6").Value = 0.0728744939271255
$ws.Range("S6").Value = 0.3643724696356275
$ws.Range("B7").Value = 0.1238532110091743
$ws.Range("D7").Value = 0.01834862385321101
$ws.Range("F7").Value = 0.03669724770642202
$ws.Range("J7").Value = 0.1422018348623853
$ws.Range("O7").Value = 0.009174311926605505
$ws.Range("Q7").Value = 0.2018348623853211
$ws.Range("R7").Value = 0.09174311926605505
$ws.Range("S7").Value = 0.3761467889908257
$ws.Range("B8").Value = 0.08370044052863436
$ws.Range("D8").Value = 0.02643171806167401
$ws.Range("F8").Value = 0.05947136563876652
$ws.Range("J8").Value = 0.1497797356828194
$ws.Range("O8").Value = 0.01101321585903084
$ws.Range("Q8").Value = 0.1629955947136564
$ws.Range("R8").Value = 0.08590308370044053
$ws.Range("S8").Value = 0.420704845814978
$ws.Range("B9").Value = 0.1232876712328767
$ws.Range("D9").Value = 0.0182648401826484
$ws.Range("F9").Value = 0.091324200913242
$ws.Range("J9").Value = 0.1141552511415525
$ws.Range("O9").Value = 0.0136986301369863
$ws.Range("Q9").Value = 0.1506849315068493
$ws.Range("R9").Value = 0.0821917808219178
$ws.Range("S9").Value = 0.4063926940639269
$ws.Range("B10").Value = 0.1172614712308813
$ws.Range("D10").Value = 0.02039329934450109
$ws.Range("F10").Value = 0.07938820101966497
$ws.Range("J10").Value = 0.1092498179169701
$ws.Range("O10").Value = 0.02039329934450109
$ws.Range("Q10").Value = 0.1886380189366351
$ws.Range("R10").Value = 0.07064821558630735
$ws.Range("S10").Value = 0.3940276766205389
$ws.Range("G11").Value = 0.1492537313432836
$ws.Range("J11").Value = 0.1074626865671642
$ws.Range("K11").Value = 0.2119402985074627
$ws.Range("L11").Value = 0.5253731343283582
$ws.Range("S11").Value = 0.005970149253731343
$ws.Range("G12").Value = 0.7204301075268817
$ws.Range("J12").Value = 0.1827956989247312
$ws.Range("K12").Value = 0.005376344086021506
$ws.Range("L12").Value = 0.04838709677419355
$ws.Range("S12").Value = 0.04301075268817205
$ws.Range("G13").Value = 0.78
$ws.Range("J13").Value = 0.22
$ws.Range("F15").Value = 0.01680672268907563
$ws.Range("H15").Value = 0.1218487394957983
$ws.Range("I15").Value = 0.07563025210084033
$ws.Range("J15").Value = 0.3403361344537815
$ws.Range("K15").Value = 0.0546218487394958
$ws.Range("M15").Value = 0.03361344537815126
$ws.Range("O15").Value = 0.03361344537815126
$ws.Range("S15").Value = 0.3235294117647059
$ws.Range("F16").Value = 0.01904761904761905
$ws.Range("H16").Value = 0.1571428571428571
$ws.Range("I16").Value = 0.06190476190476191
$ws.Range("J16").Value = 0.3904761904761905
$ws.Range("K16").Value = 0.1333333333333333
$ws.Range("O16").Value = 0.08095238095238096
$ws.Range("S16").Value = 0.1571428571428571
$ws.Range("F17").Value = 0.0248868778280543
$ws.Range("H17").Value = 0.1809954751131222
$ws.Range("I17").Value = 0.07918552036199095
$ws.Range("J17").Value = 0.4253393665158371
$ws.Range("K17").Value = 0.09049773755656108
$ws.Range("M17").Value = 0.0248868778280543
$ws.Range("O17").Value = 0.07013574660633484
$ws.Range("S17").Value = 0.1040723981900453
$ws.Range("F18").Value = 0.0154639175257732
$ws.Range("H18").Value = 0.1752577319587629
$ws.Range("I18").Value = 0.08762886597938144
$ws.Range("J18").Value = 0.3865979381443299
$ws.Range("K18").Value = 0.09278350515463918
$ws.Range("M18").Value = 0.02061855670103093
$ws.Range("O18").Value = 0.1030927835051546
$ws.Range("S18").Value = 0.1185567010309278
$ws.Range("F19").Value = 0.01480959097320169
$ws.Range("H19").Value = 0.1967559943582511
$ws.Range("I19").Value = 0.09590973201692525
$ws.Range("J19").Value = 0.3765867418899859
$ws.Range("K19").Value = 0.113540197461213
$ws.Range("M19").Value = 0.01974612129760226
$ws.Range("N19").Value = 0.002820874471086037
$ws.Range("O19").Value = 0.06064880112834979
$ws.Range("S19").Value = 0.119181946403385

Write-Host "Applied 105 cell updates"
